# Update BOM worksheet: swap the USB Logic Analyzer line item for an
# Analog Discovery 2, and adjust purchase quantities for the jack and
# jumper wires. Finally move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: was "USB Logic Analyzer - 25MHz/8-Channel" -> "Analog Discovery 2"
$ws.Range("B6").Value = "https://store.digilentinc.com/analog-discovery-2-100msps-usb-oscilloscope-logic-analyzer-and-variable-power-supply/"
$ws.Range("A6").Value = "Analog Discovery 2"
$ws.Range("C6").Value = 279
$ws.Range("D6").Value = 1

# Row 8 (3.5mm/2.5mm jack): quantity 3 -> 10
$ws.Range("D8").Value = 10

# Row 9 (jumper wires): quantity 1 -> 3
$ws.Range("D9").Value = 3

# Move the active cell selection to C15 (matches the saved view state)
$ws.Range("C15").Select() | Out-Null
